$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "selva"
$ws.Range("B2").Value = 400.0
$ws.Range("C2").Value = 910.0
$ws.Range("D2").Value = 78.0
$ws.Range("E2").Value = "coder"
$ws.Range("F2").Value = "Absorbing"
